$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1672354948805461
$ws.Range("C2").Value = 0.6040955631399317
$ws.Range("J2").Value = 0.0136518771331058
$ws.Range("P2").Value = 0.1160409556313993
$ws.Range("S2").Value = 0.09897610921501707
$ws.Range("C3").Value = 0.04371584699453552
$ws.Range("J3").Value = 0.0273224043715847
$ws.Range("P3").Value = 0.6994535519125683
$ws.Range("S3").Value = 0.2295081967213115
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.2105263157894737
$ws.Range("B6").Value = 0.07203389830508475
$ws.Range("F6").Value = 0.02966101694915254
$ws.Range("J6").Value = 0.2245762711864407
$ws.Range("O6").Value = 0.01694915254237288
$ws.Range("Q6").Value = 0.2415254237288136
$ws.Range("R6").Value = 0.07627118644067797
$ws.Range("S6").Value = 0.3389830508474576
$ws.Range("B7").Value = 0.09183673469387756
$ws.Range("D7").Value = 0.00510204081632653
$ws.Range("F7").Value = 0.03571428571428571
$ws.Range("J7").Value = 0.09693877551020408
$ws.Range("O7").Value = 0.01530612244897959
$ws.Range("Q7").Value = 0.1989795918367347
$ws.Range("R7").Value = 0.08673469387755102
$ws.Range("S7").Value = 0.4693877551020408
$ws.Range("B8").Value = 0.07377049180327869
$ws.Range("D8").Value = 0.01024590163934426
$ws.Range("F8").Value = 0.06352459016393443
$ws.Range("J8").Value = 0.110655737704918
$ws.Range("O8").Value = 0.02049180327868852
$ws.Range("Q8").Value = 0.2295081967213115
$ws.Range("R8").Value = 0.0860655737704918
$ws.Range("S8").Value = 0.4057377049180328
$ws.Range("B9").Value = 0.1050583657587549
$ws.Range("D9").Value = 0.01167315175097276
$ws.Range("F9").Value = 0.09727626459143969
$ws.Range("J9").Value = 0.1089494163424125
$ws.Range("O9").Value = 0.01167315175097276
$ws.Range("Q9").Value = 0.2334630350194553
$ws.Range("R9").Value = 0.0933852140077821
$ws.Range("S9").Value = 0.3385214007782101
$ws.Range("B10").Value = 0.1090487238979118
$ws.Range("D10").Value = 0.02242846094354215
$ws.Range("E10").Value = 0.001546790409899459
$ws.Range("F10").Value = 0.07037896365042537
$ws.Range("J10").Value = 0.08507347254447023
$ws.Range("O10").Value = 0.0185614849187935
$ws.Range("Q10").Value = 0.2382057231245166
$ws.Range("R10").Value = 0.0920340293890178
$ws.Range("S10").Value = 0.362722351121423
$ws.Range("G11").Value = 0.09507042253521127
$ws.Range("J11").Value = 0.09154929577464789
$ws.Range("K11").Value = 0.1654929577464789
$ws.Range("L11").Value = 0.6161971830985915
$ws.Range("S11").Value = 0.03169014084507042
$ws.Range("G12").Value = 0.8111111111111111
$ws.Range("J12").Value = 0.1222222222222222
$ws.Range("K12").Value = 0.01111111111111111
$ws.Range("L12").Value = 0.02222222222222222
$ws.Range("S12").Value = 0.03333333333333333
$ws.Range("G13").Value = 0.6511627906976745
$ws.Range("J13").Value = 0.3255813953488372
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.01680672268907563
$ws.Range("H15").Value = 0.1764705882352941
$ws.Range("I15").Value = 0.05882352941176471
$ws.Range("J15").Value = 0.3529411764705883
$ws.Range("K15").Value = 0.06302521008403361
$ws.Range("M15").Value = 0.01680672268907563
$ws.Range("N15").Value = 0.004201680672268907
$ws.Range("O15").Value = 0.06722689075630252
$ws.Range("S15").Value = 0.2436974789915966
$ws.Range("F16").Value = 0.03825136612021858
$ws.Range("H16").Value = 0.2185792349726776
$ws.Range("I16").Value = 0.1038251366120219
$ws.Range("J16").Value = 0.3497267759562842
$ws.Range("K16").Value = 0.1147540983606557
$ws.Range("M16").Value = 0.01639344262295082
$ws.Range("O16").Value = 0.0273224043715847
$ws.Range("S16").Value = 0.1311475409836066
$ws.Range("F17").Value = 0.01188455008488964
$ws.Range("H17").Value = 0.166383701188455
$ws.Range("I17").Value = 0.1086587436332767
$ws.Range("J17").Value = 0.4125636672325976
$ws.Range("K17").Value = 0.08488964346349745
$ws.Range("M17").Value = 0.01528013582342954
$ws.Range("N17").Value = 0.001697792869269949
$ws.Range("O17").Value = 0.06960950764006792
$ws.Range("S17").Value = 0.1290322580645161
$ws.Range("F18").Value = 0.00904977375565611
$ws.Range("H18").Value = 0.2036199095022624
$ws.Range("I18").Value = 0.1402714932126697
$ws.Range("J18").Value = 0.4343891402714932
$ws.Range("K18").Value = 0.09049773755656108
$ws.Range("M18").Value = 0.01357466063348416
$ws.Range("O18").Value = 0.05429864253393665
$ws.Range("S18").Value = 0.05429864253393665
$ws.Range("F19").Value = 0.02530674846625767
$ws.Range("H19").Value = 0.200920245398773
$ws.Range("I19").Value = 0.09892638036809816
$ws.Range("J19").Value = 0.3711656441717792
$ws.Range("K19").Value = 0.09585889570552147
$ws.Range("M19").Value = 0.01993865030674847
$ws.Range("N19").Value = 0.0007668711656441718
$ws.Range("O19").Value = 0.07131901840490798
$ws.Range("S19").Value = 0.1157975460122699
